# ROBE-186 import error for big decimal numbers larger than 11 digits
# Bump the sample "F" column values past 11 digits so the sample workbook
# exercises the big-decimal import path, widen column F so the bigger
# numbers are fully visible, and leave the selection on the column that
# was touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 11111111111
$ws.Range("F2").Value = 11111111112
$ws.Range("F3").Value = 11111111113
$ws.Range("F4").Value = 11111111114
$ws.Range("F5").Value = 11111111115

# Widen column F so the 11-digit values fit.
$ws.Columns.Item(6).ColumnWidth = 13.5

# Mirror the selection recorded in the edited workbook.
$ws.Range("F1:F5").Select()
